# Update selected balance sheet figures on the "PWR" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWR")

# Row 4 - Inventory
$ws.Range("B4").Value = 50000000.0
$ws.Range("C4").Value = 40000000.0
$ws.Range("D4").Value = 48000000.0
$ws.Range("E4").Value = 50000000.0
$ws.Range("F4").Value = 56000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 798000000.0
$ws.Range("C14").Value = 1544000000.0
$ws.Range("D14").Value = 1289000000.0
$ws.Range("E14").Value = 1360000000.0
$ws.Range("F14").Value = 1490000000.0

# Row 19 - Long Term Tax Liability (Deferred)
$ws.Range("B19").Value = 149000000.0
$ws.Range("C19").Value = 206000000.0
$ws.Range("D19").Value = 221000000.0
$ws.Range("E19").Value = 218000000.0
$ws.Range("F19").Value = 215000000.0
